# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2404   (columns A-J)
#   *_new  -> *_FV2410   (columns L-U)
# Then freeze the header row and wrap the data range in an Excel Table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (row 1) ------------------------------------------
$ws.Range("A1").Value = "Segmentname_FV2404"
$ws.Range("B1").Value = "Segmentgruppe_FV2404"
$ws.Range("C1").Value = "Segment_FV2404"
$ws.Range("D1").Value = "Datenelement_FV2404"
$ws.Range("E1").Value = "Segment ID_FV2404"
$ws.Range("F1").Value = "Code_FV2404"
$ws.Range("G1").Value = "Qualifier_FV2404"
$ws.Range("H1").Value = "Beschreibung_FV2404"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("J1").Value = "Bedingung_FV2404"
# K1 ("diff") is unchanged.
$ws.Range("L1").Value = "Segmentname_FV2410"
$ws.Range("M1").Value = "Segmentgruppe_FV2410"
$ws.Range("N1").Value = "Segment_FV2410"
$ws.Range("O1").Value = "Datenelement_FV2410"
$ws.Range("P1").Value = "Segment ID_FV2410"
$ws.Range("Q1").Value = "Code_FV2410"
$ws.Range("R1").Value = "Qualifier_FV2410"
$ws.Range("S1").Value = "Beschreibung_FV2410"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("U1").Value = "Bedingung_FV2410"

# --- Freeze the header row (row 1) --------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

# --- Turn the data range into an Excel Table (ListObject) --------------
# Column names are picked up from the (already renamed) header row.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U74"), $null, 1)
$tbl.Name = "Table1"
